# Apply: model and template with unit, description and enum
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers ---
# Column H/I headers swap due to shared-string reindexing (SamplePortion now
# precedes Result in the shared string table, and the H/I cells point at the
# same underlying string indices as before, which now resolve differently).
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# --- Row 2: types, now annotated with units ---
$ws.Range("A2").Value = "#string"
$ws.Range("B2").Value = "#string"
$ws.Range("C2").Value = "#date"
$ws.Range("D2").Value = "#string"
$ws.Range("E2").Value = "#string"
$ws.Range("F2").Value = "#string"
$ws.Range("G2").Value = "#string"
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"
$ws.Range("J2").Value = "#float,  unit:celsius"
$ws.Range("K2").Value = "#integer,  unit:hours"
$ws.Range("L2").Value = "#string"

# --- Row 3: new description / enum row ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#Temperature"
$ws.Range("K3").Value = "#Temps"
$ws.Range("L3").Value = "#CycleDeTemperature"
